# Auto-generated script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.769.27"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "3.253.13"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'578.75"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "'181.76"
$ws.Range("E6").Value = "  +4.74%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").Value = "'0.133"
$ws.Range("E9").Value = "  +4.62%  "
$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").Value = "'0.415"
$ws.Range("E11").Value = "  +4.29%  "
$ws.Range("D12").Value = "3.819.97"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "'28.58"
$ws.Range("E14").Value = "  +4.23%  "
$ws.Range("D15").Value = "67.768.02"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").Value = "'0.0000172"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("D17").Value = "3.261.45"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").Value = "'5.83"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "'13.51"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").Value = "'379.08"
$ws.Range("E20").Value = "  +3.81%  "
$ws.Range("D21").Value = "'7.65"
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'71.40"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D24").Value = "'0.512"
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("D25").Value = "'0.0000119"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("D26").Value = "'9.92"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +2.22%  "
$ws.Range("D30").Value = "'5.66"
$ws.Range("E30").Value = "  +3.40%  "
$ws.Range("E31").Value = "  +3.21%  "
$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'7.00"
$ws.Range("E33").Value = "  +5.02%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.27"
$ws.Range("E34").Value = "  +5.09%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +5.63%  "
$ws.Range("D36").Value = "'163.62"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'1.88"
$ws.Range("E37").Value = "  +3.70%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Value = "'0.847"
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'4.64"
$ws.Range("E39").Value = "  +10.06%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'26.54"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'6.64"
$ws.Range("E41").Value = "  +6.39%  "
$ws.Range("D42").Value = "'2.60"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'348.55"
$ws.Range("E43").Value = "  +5.31%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'25.57"
$ws.Range("E44").Value = "  +5.04%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'41.04"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0683"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.615.36"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").Value = "'0.994"
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.17"
$ws.Range("E51").Value = "  +3.52%  "

Write-Host "Applied 116 cell updates"
